# boot strapping + correlation plots HFT
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "exclude" flags (column L) on existing rows ---
$ws.Range("L6").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("L17").Value = 1
$ws.Range("L19").Value = 0
$ws.Range("L20").Value = 1

# --- Replace literal EF value on row 16 with an average formula ---
$ws.Range("F16").Formula = "=(0.0127+0.0133)/2"

# --- Append three new data rows (21-23) ---
$ws.Range("A21").Value = "2024-07-18_C_e"
$ws.Range("B21").Value = "freq"
$ws.Range("C21").Value = 0.82
$ws.Range("D21").Value = 0.0002
$ws.Range("E21").Value = 0.15
$ws.Range("F21").Value = 0.0131
$ws.Range("G21").Value = 47
$ws.Range("H21").Value = 202.1
$ws.Range("I21").Value = 47.2159
$ws.Range("J21").Value = "Blackman"
$ws.Range("L21").Value = 0

$ws.Range("A22").Value = "2024-07-18_D_e"
$ws.Range("B22").Value = "freq"
$ws.Range("C22").Value = 0.82
$ws.Range("D22").Value = 0.0002
$ws.Range("E22").Value = 0.2
$ws.Range("F22").Value = 0.0131
$ws.Range("G22").Value = 47
$ws.Range("H22").Value = 202.1
$ws.Range("I22").Value = 47.2159
$ws.Range("J22").Value = "Blackman"
$ws.Range("L22").Value = 0

$ws.Range("A23").Value = "2024-07-18_E_e"
$ws.Range("B23").Value = "freq"
$ws.Range("C23").Value = 0.82
$ws.Range("D23").Value = 0.0002
$ws.Range("E23").Value = 0.2
$ws.Range("F23").Value = 0.0131
$ws.Range("G23").Value = 47
$ws.Range("H23").Value = 202.1
$ws.Range("I23").Value = 47.2159
$ws.Range("J23").Value = "Blackman"
$ws.Range("L23").Value = 0

# --- Update view: move the active selection (also clears the stale topLeftCell scroll state) ---
$ws.Range("O9").Select()
